# Generate Report for handback
#
# For each localized-language sheet (zh-cn, de-de) the localization files
# listed in rows 2-3 have now been handed back: the "Status" column moves
# from "Ready for handoff" to "Handed back: in sync with en-us", the
# "Latest Target File" (E) / "Latest Handback File" (F) columns are filled
# in (mirroring the source .md file and the handed-off .xlf file), and the
# "Latest Handback DateTime" (G) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$handbackStatus = "Handed back: in sync with en-us"

$sheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-01-08 18:58:49" },
    @{ Name = "de-de"; HandbackTime = "2016-01-08 18:59:05" }
)

# BGR-encoded equivalent of the workbook's custom hyperlink colour (RGB FF6495ED)
$hyperlinkColor = 15570276

# The "Overview" sheet mirrors the same Status text per language (columns B/C
# for rows 2/3), sharing the same underlying string -- update it too so it
# stays in sync with the detail sheets below.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $handbackStatus
$wsOverview.Range("C2").Value = $handbackStatus
$wsOverview.Range("B3").Value = $handbackStatus
$wsOverview.Range("C3").Value = $handbackStatus

foreach ($sheetInfo in $sheets) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    foreach ($row in 2, 3) {

        # Status column: handoff -> handback
        $ws.Cells.Item($row, 2).Value = $handbackStatus

        $sourceDisplay = $ws.Cells.Item($row, 1).Value2
        $sourceAddress = $ws.Hyperlinks.Item(1).Address
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Row -eq $row -and $h.Range.Column -eq 1) {
                $sourceAddress = $h.Address
            }
        }

        $handoffDisplay = $ws.Cells.Item($row, 3).Value2
        $handoffAddress = $sourceAddress
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Row -eq $row -and $h.Range.Column -eq 3) {
                $handoffAddress = $h.Address
            }
        }

        # E = Latest Target File (same file as the source markdown, column A)
        $eCell = $ws.Cells.Item($row, 5)
        $ws.Hyperlinks.Add($eCell, $sourceAddress, "", "", $sourceDisplay)
        $eCell.Font.Underline = $true
        $eCell.Font.Color = $hyperlinkColor

        # F = Latest Handback File (the .xlf that was handed off, column C)
        $fCell = $ws.Cells.Item($row, 6)
        $ws.Hyperlinks.Add($fCell, $handoffAddress, "", "", $handoffDisplay)
        $fCell.Font.Underline = $true
        $fCell.Font.Color = $hyperlinkColor

        # G = Latest Handback DateTime
        $ws.Cells.Item($row, 7).Value = $sheetInfo.HandbackTime
    }
}
